$wb = $excel.ActiveWorkbook

# zh-cn sheet: row for e86a555d-... handoff/handback pair
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-23 02:43:35"
$wsZh.Range("H3").Value = "2016-03-23 02:44:01"

# de-de sheet: row for e86a555d-... handoff/handback pair
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-23 02:43:39"
$wsDe.Range("H3").Value = "2016-03-23 02:44:08"
